# Updated cryptos list on Sun Jul 23 04:35:47 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures in the crypto table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.846.46'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = '1.869.67'
$ws.Range("E3").Value = '  -1.40%  '
$ws.Range("D4").Value = '''0.9990'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '''0.7401'
$ws.Range("E5").Value = '  -4.59%  '
$ws.Range("D6").Value = '''241.94'
$ws.Range("E6").Value = '  -1.18%  '
$ws.Range("D7").Value = '''0.9996'
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = '''0.3153'
$ws.Range("E8").Value = '  +0.30%  '
$ws.Range("D9").Value = '''24.67'
$ws.Range("E9").Value = '  -4.43%  '
$ws.Range("D10").Value = '''0.07105'
$ws.Range("E10").Value = '  -2.28%  '
$ws.Range("D11").Value = '''0.08373'
$ws.Range("E11").Value = '  -6.29%  '
$ws.Range("D12").Value = '''0.7520'
$ws.Range("E12").Value = '  -3.06%  '
$ws.Range("D13").Value = '''5.448'
$ws.Range("E13").Value = '  -0.27%  '
$ws.Range("D14").Value = '1.866.73'
$ws.Range("E14").Value = '  +0.07%  '
$ws.Range("D15").Value = '''92.49'
$ws.Range("E15").Value = '  -2.52%  '
$ws.Range("D16").Value = '29.853.91'
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("D17").Value = '''6.028'
$ws.Range("E17").Value = '  -2.84%  '
$ws.Range("E18").Value = '  -3.16%  '
$ws.Range("D19").Value = '''242.84'
$ws.Range("E19").Value = '  -1.67%  '
$ws.Range("D20").Value = '''0.000007817'
$ws.Range("E20").Value = '  -1.23%  '
$ws.Range("D21").Value = '''0.9984'
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("D22").Value = '2.118.38'
$ws.Range("E22").Value = '  -0.58%  '
$ws.Range("D23").Value = '''7.925'
$ws.Range("E23").Value = '  -2.86%  '
$ws.Range("D24").Value = '''0.9994'
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").Value = '''0.1565'
$ws.Range("E25").Value = '  -1.68%  '
$ws.Range("D26").Value = '''9.293'
$ws.Range("E26").Value = '  -2.82%  '
$ws.Range("D27").Value = '''164.10'
$ws.Range("E27").Value = '  +0.55%  '
$ws.Range("E28").Value = '  -1.59%  '
$ws.Range("D29").Value = '''2.018'
$ws.Range("E29").Value = '  -1.70%  '
$ws.Range("D30").Value = '''1.472'
$ws.Range("E30").Value = '  +3.11%  '
$ws.Range("D31").Value = '''4.634'
$ws.Range("E31").Value = '  +2.16%  '
$ws.Range("D32").Value = '''1.534'
$ws.Range("E32").Value = '  -0.96%  '
$ws.Range("D33").Value = '''4.305'
$ws.Range("E33").Value = '  +4.34%  '
$ws.Range("D34").Value = '''0.05313'
$ws.Range("E34").Value = '  -4.01%  '
$ws.Range("E35").Value = '  -1.42%  '
$ws.Range("D36").Value = '''0.7520'
$ws.Range("E36").Value = '  -0.47%  '
$ws.Range("D37").Value = '''0.9997'
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("D38").Value = '''2.698'
$ws.Range("E38").Value = '  -0.80%  '
$ws.Range("E39").Value = '  -0.87%  '
$ws.Range("D40").Value = '''2.749'
$ws.Range("D41").Value = '''0.4472'
$ws.Range("E41").Value = '  -1.26%  '
$ws.Range("D42").Value = '1.102.57'
$ws.Range("E42").Value = '  +1.47%  '
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("D44").Value = '''72.14'
$ws.Range("E44").Value = '  -2.73%  '
$ws.Range("D45").Value = '''0.8593'
$ws.Range("E45").Value = '  +0.41%  '
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").Value = '''103.09'
$ws.Range("E47").Value = '  +0.24%  '
$ws.Range("E48").Value = '  +0.63%  '
$ws.Range("D49").Value = '''1.839'
$ws.Range("E49").Value = '  -3.13%  '
$ws.Range("D50").Value = '''3.046'
$ws.Range("E50").Value = '  +1.54%  '
$ws.Range("D51").Value = '2.014.55'
$ws.Range("E51").Value = '  -0.40%  '
